$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.586.33"
$ws.Range("D3").Value = "2.087.95"
$ws.Range("E3").Value = "  -1.40%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.010"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "341.97"
$ws.Range("E5").Value = "  -2.39%  "
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5122"
$ws.Range("E7").Value = "  -2.69%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4373"
$ws.Range("E8").Value = "  -3.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.25"
$ws.Range("E9").Value = "  -2.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09089"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.64"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "2.098.28"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.733"
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.104"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "100.19"
$ws.Range("E16").Value = "  -1.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.009"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001141"
$ws.Range("E18").Value = "  -2.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "20.93"
$ws.Range("E19").Value = "  +7.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06640"
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.009"
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.151"
$ws.Range("E22").Value = "  -2.65%  "
$ws.Range("D23").Value = "29.606.15"
$ws.Range("E23").Value = "  -3.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.56"
$ws.Range("E24").Value = "  -2.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.302"
$ws.Range("E25").Value = "  -3.58%  "
$ws.Range("D26").Value = "2.330.20"
$ws.Range("E26").Value = "  -1.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.77"
$ws.Range("E27").Value = "  -2.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.01"
$ws.Range("E28").Value = "  -1.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.508"
$ws.Range("E29").Value = "  -1.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.04"
$ws.Range("E30").Value = "  -3.70%  "
$ws.Range("E31").Value = "  -6.20%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1042"
$ws.Range("E32").Value = "  -3.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.617"
$ws.Range("E33").Value = "  -2.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.110"
$ws.Range("E34").Value = "  -4.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.961"
$ws.Range("E35").Value = "  -1.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.988"
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.15"
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02565"
$ws.Range("E38").Value = "  -3.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06654"
$ws.Range("E39").Value = "  -3.03%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2224"
$ws.Range("E40").Value = "  -4.30%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6825"
$ws.Range("E41").Value = "  -1.22%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.31"
$ws.Range("E42").Value = "  -2.33%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.278"
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6684"
$ws.Range("E44").Value = "  +3.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.11"
$ws.Range("E45").Value = "  -4.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.277"
$ws.Range("E46").Value = "  -2.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.601"
$ws.Range("E47").Value = "  -4.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.214"
$ws.Range("E48").Value = "  -3.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000336"
$ws.Range("E49").Value = "  -5.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "81.31"
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.173"
$ws.Range("E51").Value = "  -1.67%  "
